# Auto-generated from diff: set literal numeric cell values to match
# the refreshed TPM (transcripts-per-million) recompute for the
# Fstl1 -> Dip2a ligand-receptor sheet (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 6.794730333333334
$ws.Range("H2").Value2 = 20.384191
$ws.Range("I2").Value2 = 0.01549131991191145
$ws.Range("J2").Value2 = 0.01549131991191145
$ws.Range("M2").Value2 = 6.579711666666667
$ws.Range("N2").Value2 = 19.739135
$ws.Range("O2").Value2 = 0.3046157543678319
$ws.Range("P2").Value2 = 0.3046157543678319
$ws.Range("Q2").Value2 = 44.70736644608723
$ws.Range("R2").Value2 = 402.3662980147851
$ws.Range("S2").Value2 = 0.004718900101120324
$ws.Range("T2").Value2 = 0.004718900101120323
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 6.794730333333334
$ws.Range("H3").Value2 = 20.384191
$ws.Range("I3").Value2 = 0.01549131991191145
$ws.Range("J3").Value2 = 0.01549131991191145
$ws.Range("M3").Value2 = 9.033654666666669
$ws.Range("O3").Value2 = 0.4182240302300713
$ws.Range("P3").Value2 = 0.4182240302300712
$ws.Range("Q3").Value2 = 61.38124738445824
$ws.Range("R3").Value2 = 552.4312264601241
$ws.Range("S3").Value2 = 0.006478842247142961
$ws.Range("T3").Value2 = 0.00647884224714296
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 6.794730333333334
$ws.Range("H4").Value2 = 20.384191
$ws.Range("I4").Value2 = 0.01549131991191145
$ws.Range("J4").Value2 = 0.01549131991191145
$ws.Range("M4").Value2 = 5.986671
$ws.Range("N4").Value2 = 17.960013
$ws.Range("O4").Value2 = 0.2771602154020968
$ws.Range("P4").Value2 = 0.2771602154020968
$ws.Range("Q4").Value2 = 40.67781503938701
$ws.Range("R4").Value2 = 366.100335354483
$ws.Range("S4").Value2 = 0.00429357756364817
$ws.Range("T4").Value2 = 0.00429357756364817
$ws.Range("G5").Value2 = 398.9908546666666
$ws.Range("I5").Value2 = 0.9096600848522713
$ws.Range("J5").Value2 = 0.9096600848522712
$ws.Range("M5").Value2 = 6.579711666666667
$ws.Range("N5").Value2 = 19.739135
$ws.Range("O5").Value2 = 0.3046157543678319
$ws.Range("P5").Value2 = 0.3046157543678319
$ws.Range("Q5").Value2 = 2625.244781343571
$ws.Range("R5").Value2 = 23627.20303209214
$ws.Range("S5").Value2 = 0.2770967929655807
$ws.Range("T5").Value2 = 0.2770967929655806
$ws.Range("G6").Value2 = 398.9908546666666
$ws.Range("I6").Value2 = 0.9096600848522713
$ws.Range("J6").Value2 = 0.9096600848522712
$ws.Range("M6").Value2 = 9.033654666666669
$ws.Range("O6").Value2 = 0.4182240302300713
$ws.Range("P6").Value2 = 0.4182240302300712
$ws.Range("S6").Value2 = 0.3804417068263455
$ws.Range("T6").Value2 = 0.3804417068263454
$ws.Range("G7").Value2 = 398.9908546666666
$ws.Range("I7").Value2 = 0.9096600848522713
$ws.Range("J7").Value2 = 0.9096600848522712
$ws.Range("M7").Value2 = 5.986671
$ws.Range("N7").Value2 = 17.960013
$ws.Range("O7").Value2 = 0.2771602154020968
$ws.Range("P7").Value2 = 0.2771602154020968
$ws.Range("Q7").Value2 = 2388.626978898148
$ws.Range("R7").Value2 = 21497.64281008333
$ws.Range("S7").Value2 = 0.2521215850603452
$ws.Range("T7").Value2 = 0.2521215850603452
$ws.Range("G8").Value2 = 32.82974100000001
$ws.Range("H8").Value2 = 98.48922300000001
$ws.Range("I8").Value2 = 0.0748485952358172
$ws.Range("J8").Value2 = 0.07484859523581719
$ws.Range("M8").Value2 = 6.579711666666667
$ws.Range("N8").Value2 = 19.739135
$ws.Range("O8").Value2 = 0.3046157543678319
$ws.Range("P8").Value2 = 0.3046157543678319
$ws.Range("Q8").Value2 = 216.010229871345
$ws.Range("R8").Value2 = 1944.092068842105
$ws.Range("S8").Value2 = 0.02280006130113097
$ws.Range("T8").Value2 = 0.02280006130113096
$ws.Range("G9").Value2 = 32.82974100000001
$ws.Range("H9").Value2 = 98.48922300000001
$ws.Range("I9").Value2 = 0.0748485952358172
$ws.Range("J9").Value2 = 0.07484859523581719
$ws.Range("M9").Value2 = 9.033654666666669
$ws.Range("O9").Value2 = 0.4182240302300713
$ws.Range("P9").Value2 = 0.4182240302300712
$ws.Range("Q9").Value2 = 296.5725429901081
$ws.Range("R9").Value2 = 2669.152886910973
$ws.Range("S9").Value2 = 0.03130348115658278
$ws.Range("T9").Value2 = 0.03130348115658277
$ws.Range("G10").Value2 = 32.82974100000001
$ws.Range("H10").Value2 = 98.48922300000001
$ws.Range("I10").Value2 = 0.0748485952358172
$ws.Range("J10").Value2 = 0.07484859523581719
$ws.Range("M10").Value2 = 5.986671
$ws.Range("N10").Value2 = 17.960013
$ws.Range("O10").Value2 = 0.2771602154020968
$ws.Range("P10").Value2 = 0.2771602154020968
$ws.Range("Q10").Value2 = 196.540858382211
$ws.Range("R10").Value2 = 1768.867725439899
$ws.Range("S10").Value2 = 0.02074505277810345
$ws.Range("T10").Value2 = 0.02074505277810345
Write-Output "Updated 106 cells with refreshed TPM values."
